# Updates cryptos list data (price/volume columns + row 20/21 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.607.44"
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("D3").Value = "2.970.26"
$ws.Range("E3").Value = "  -5.15%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.64"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.90"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.67%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").Value = "2.979.88"
$ws.Range("E9").Value = "  -5.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.114"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.15"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.05%  "
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").Value = "3.488.52"
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.125"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "61.677.09"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.81"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").Value = "2.974.29"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.03"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.40"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.71"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.88"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").Value = "3.092.83"
$ws.Range("E26").Value = "  -5.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.188"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "0.0₃0942"
$ws.Range("E29").Value = "  -6.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.31"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.73"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.46"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "160.69"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("E36").Value = "  -3.81%  "
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("D41").Value = "2.414.38"
$ws.Range("E41").Value = "  -9.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.28"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.21"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.16%  "
$ws.Range("E44").Value = "  -3.13%  "
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("E46").Value = "  -1.92%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.04"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.93"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "269.94"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0954"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.68%  "
